$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.774.96"
$ws.Range("E2").Value = "  -1.68%  "

$ws.Range("D3").Value = "2.351.24"
$ws.Range("E3").Value = "  -2.63%  "

$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").Value = "'321.03"
$ws.Range("E5").Value = "  -0.90%  "

$ws.Range("E6").Value = "  +1.13%  "

$ws.Range("D7").Value = "'0.638"
$ws.Range("E7").Value = "  -1.04%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "'0.617"
$ws.Range("E9").Value = "  -6.77%  "

$ws.Range("D10").Value = "'41.05"
$ws.Range("E10").Value = "  -2.69%  "

$ws.Range("D11").Value = "'0.0925"
$ws.Range("E11").Value = "  -2.69%  "

$ws.Range("D12").Value = "'8.44"
$ws.Range("E12").Value = "  -2.49%  "

$ws.Range("D13").Value = "'0.995"
$ws.Range("E13").Value = "  -4.91%  "

$ws.Range("E14").Value = "  -0.15%  "

$ws.Range("D15").Value = "'16.02"
$ws.Range("E15").Value = "  -7.42%  "

$ws.Range("D16").Value = "2.704.00"
$ws.Range("E16").Value = "  -2.81%  "

$ws.Range("D17").Value = "2.347.21"
$ws.Range("E17").Value = "  -2.94%  "

$ws.Range("D18").Value = "42.719.72"
$ws.Range("E18").Value = "  -1.92%  "

$ws.Range("D19").Value = "'7.72"
$ws.Range("E19").Value = "  +3.34%  "

$ws.Range("E20").Value = "  -3.68%  "

$ws.Range("D21").Value = "'77.28"
$ws.Range("E21").Value = "  +2.32%  "

$ws.Range("D22").Value = "'3.58"
$ws.Range("E22").Value = "  +2.38%  "

$ws.Range("D23").Value = "'260.25"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").Value = "'2.33"
$ws.Range("E24").Value = "  -5.15%  "

$ws.Range("D25").Value = "'9.59"
$ws.Range("E25").Value = "  -1.03%  "

$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("D27").Value = "'11.44"
$ws.Range("E27").Value = "  -4.16%  "

$ws.Range("D28").Value = "'23.30"
$ws.Range("E28").Value = "  +2.16%  "

$ws.Range("E29").Value = "  -0.83%  "

$ws.Range("D30").Value = "'174.77"
$ws.Range("E30").Value = "  -2.65%  "

$ws.Range("D31").Value = "'36.24"
$ws.Range("E31").Value = "  -4.89%  "

$ws.Range("D32").Value = "'3.00"
$ws.Range("E32").Value = "  -7.06%  "

$ws.Range("D33").Value = "'0.0892"
$ws.Range("E33").Value = "  -4.69%  "

$ws.Range("D34").Value = "'6.09"
$ws.Range("E34").Value = "  +2.34%  "

$ws.Range("D35").Value = "'0.131"
$ws.Range("E35").Value = "  -1.44%  "

$ws.Range("D36").Value = "'0.113"
$ws.Range("E36").Value = "  +6.04%  "

$ws.Range("D37").Value = "'4.63"
$ws.Range("E37").Value = "  -5.43%  "

$ws.Range("D38").Value = "'0.0358"
$ws.Range("E38").Value = "  -3.87%  "

$ws.Range("D39").Value = "'3.81"
$ws.Range("E39").Value = "  -4.08%  "

$ws.Range("D40").Value = "'2.68"
$ws.Range("E40").Value = "  -7.59%  "

$ws.Range("D41").Value = "'1.47"
$ws.Range("E41").Value = "  -9.71%  "

$ws.Range("D42").Value = "'70.84"
$ws.Range("E42").Value = "  +1.51%  "

$ws.Range("D43").Value = "'0.232"
$ws.Range("E43").Value = "  -0.71%  "

$ws.Range("E44").Value = "  -0.22%  "

$ws.Range("D45").Value = "'114.81"
$ws.Range("E45").Value = "  -9.34%  "

$ws.Range("D46").Value = "'11.88"
$ws.Range("E46").Value = "  -5.67%  "

$ws.Range("D47").Value = "'5.51"
$ws.Range("E47").Value = "  -3.09%  "

$ws.Range("D48").Value = "'9.21"
$ws.Range("E48").Value = "  -3.98%  "

$ws.Range("D49").Value = "'83.74"
$ws.Range("E49").Value = "  +6.58%  "

$ws.Range("D50").Value = "'73.19"
$ws.Range("E50").Value = "  -0.88%  "

$ws.Range("E51").Value = "  -1.23%  "
